# Shift every transaction timestamp in column C (rows 2-51) forward by
# 15 hours 16 minutes, keeping the "yyyy-MM-dd HH:mm" text format used
# by the sheet (the cells hold plain text, not real Excel date values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $addr = "C$r"
    $raw = $ws.Range($addr).Value2

    if ([string]::IsNullOrEmpty($raw)) {
        continue
    }

    $parts = $raw.Split(" ")
    $dateParts = $parts[0].Split("-")
    $timeParts = $parts[1].Split(":")

    $y  = [int]$dateParts[0]
    $m  = [int]$dateParts[1]
    $d  = [int]$dateParts[2]
    $h  = [int]$timeParts[0]
    $mi = [int]$timeParts[1]

    $dt = Get-Date -Year $y -Month $m -Day $d -Hour $h -Minute $mi -Second 0
    $dt = $dt.AddHours(15).AddMinutes(16)

    $ws.Range($addr).Value = $dt.ToString("yyyy-MM-dd HH:mm")
}
